$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "29.711.05"
$ws.Range("E2").Value = "  +5.27%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.922.74"
$ws.Range("E3").Value = "  +3.64%  "

# Row 4 (TetherUSD)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.41"
$ws.Range("E5").Value = "  +1.79%  "

# Row 6 (USDC)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7 (XRP)
$ws.Range("E7").Value = "  +2.94%  "

# Row 8 (Cardano)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4120"
$ws.Range("E8").Value = "  +5.07%  "

# Row 9 (OKB)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.14"
$ws.Range("E9").Value = "  +1.38%  "

# Row 10 (Dogecoin)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08040"
$ws.Range("E10").Value = "  +3.17%  "

# Row 11 (Polygon)
$ws.Range("E11").Value = "  +3.84%  "

# Row 12 (Solana)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.49"
$ws.Range("E12").Value = "  +5.91%  "

# Row 13 (WrappedEther)
$ws.Range("D13").Value = "1.930.62"
$ws.Range("E13").Value = "  +2.71%  "

# Row 14 (Polkadot)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.999"
$ws.Range("E14").Value = "  +3.62%  "

# Row 15 (Chainlink)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.198"
$ws.Range("E15").Value = "  +3.69%  "

# Row 16 (Litecoin)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.27"
$ws.Range("E16").Value = "  +2.80%  "

# Row 17 (BinanceUSD)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.23%  "

# Row 18 (ShibaInu)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001034"
$ws.Range("E18").Value = "  +1.88%  "

# Row 19 (TRON)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06593"
$ws.Range("E19").Value = "  +1.09%  "

# Row 20 (Avalanche)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.85"
$ws.Range("E20").Value = "  +4.82%  "

# Row 21 (Dai)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 (WrappedBTC)
$ws.Range("D22").Value = "29.664.97"
$ws.Range("E22").Value = "  +5.13%  "

# Row 23 (Uniswap)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.582"
$ws.Range("E23").Value = "  +5.55%  "

# Row 24 (Cosmos)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.70"
$ws.Range("E24").Value = "  +9.58%  "

# Row 25 (Toncoin)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.212"
$ws.Range("E25").Value = "  -1.92%  "

# Row 26 (WrappedliquidstakedEther2.0)
$ws.Range("D26").Value = "2.159.67"
$ws.Range("E26").Value = "  +3.70%  "

# Row 27 (Monero)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.21"
$ws.Range("E27").Value = "  -0.45%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.89"
$ws.Range("E28").Value = "  +3.83%  "

# Row 29 (LidoDAOToken)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.150"
$ws.Range("E29").Value = "  +5.35%  "

# Row 30 (InternetComputer(DFINITY))
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.723"
$ws.Range("E30").Value = "  +8.90%  "

# Row 31 (BitcoinCash)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.68"
$ws.Range("E31").Value = "  +1.35%  "

# Row 32 (ImmutableX)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.066"
$ws.Range("E32").Value = "  +13.57%  "

# Row 33 (Stellar)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09490"
$ws.Range("E33").Value = "  +2.51%  "

# Row 34 (ARBITRUM)
$ws.Range("E34").Value = "  +4.95%  "

# Row 35 (HuobiToken)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.571"
$ws.Range("E35").Value = "  -0.92%  "

# Row 36 (Filecoin)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.422"
$ws.Range("E36").Value = "  +4.65%  "

# Row 37: 'VeChain' -> 'Hedera'
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06143"
$ws.Range("E37").Value = "  +2.33%  "

# Row 38: 'Hedera' -> 'VeChain'
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02275"
$ws.Range("E38").Value = "  +3.80%  "

# Row 39 (FraxShare)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.437"
$ws.Range("E39").Value = "  +3.29%  "

# Row 40 (TrustWalletToken)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("E40").Value = "  +1.95%  "

# Row 41 (TheSandbox)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5900"
$ws.Range("E41").Value = "  +4.42%  "

# Row 42 (Algorand)
$ws.Range("E42").Value = "  +3.39%  "

# Row 43 (Aptos)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.24"
$ws.Range("E43").Value = "  +2.89%  "

# Row 44 (WEMIXTOKEN)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.261"
$ws.Range("E44").Value = "  +1.24%  "

# Row 45 (RenderToken)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.357"
$ws.Range("E45").Value = "  +3.13%  "

# Row 46: 'Cronos' -> 'EnergySwap'
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.31"
$ws.Range("E46").Value = "  +4.43%  "

# Row 47: 'Decentraland' -> 'Cronos'
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07505"
$ws.Range("E47").Value = "  +4.64%  "

# Row 48: 'EnergySwap' -> 'Decentraland'
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5590"
$ws.Range("E48").Value = "  +4.44%  "

# Row 49 (NEARProtocol)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.938"
$ws.Range("E49").Value = "  +4.28%  "

# Row 50 (Quant)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.28"
$ws.Range("E50").Value = "  +3.49%  "

# Row 51 (WOONetwork)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3013"
$ws.Range("E51").Value = "  +15.22%  "
